$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Itga2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.45491533333333
$ws.Range("H2").Value = 31.364746
$ws.Range("I2").Value = 0.0134573334963438
$ws.Range("J2").Value = 0.0134573334963438
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.888921
$ws.Range("N2").Value = 5.666763
$ws.Range("O2").Value = 0.4551532417350329
$ws.Range("P2").Value = 0.4551532417350328
$ws.Range("Q2").Value = 19.74850912635533
$ws.Range("R2").Value = 177.736582137198
$ws.Range("S2").Value = 0.006125148965970325
$ws.Range("T2").Value = 0.006125148965970323

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Itga2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.45491533333333
$ws.Range("H3").Value = 31.364746
$ws.Range("I3").Value = 0.0134573334963438
$ws.Range("J3").Value = 0.0134573334963438
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.087098333333333
$ws.Range("N3").Value = 3.261295
$ws.Range("O3").Value = 0.2619465454094788
$ws.Range("P3").Value = 0.2619465454094787
$ws.Range("Q3").Value = 11.36552103400778
$ws.Range("R3").Value = 102.28968930607
$ws.Range("S3").Value = 0.003525102019790521
$ws.Range("T3").Value = 0.00352510201979052

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Itga2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.45491533333333
$ws.Range("H4").Value = 31.364746
$ws.Range("I4").Value = 0.0134573334963438
$ws.Range("J4").Value = 0.0134573334963438
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.174057666666667
$ws.Range("N4").Value = 3.522173
$ws.Range("O4").Value = 0.2829002128554884
$ws.Range("P4").Value = 0.2829002128554884
$ws.Range("Q4").Value = 12.27467350145089
$ws.Range("R4").Value = 110.472061513058
$ws.Range("S4").Value = 0.003807082510582955
$ws.Range("T4").Value = 0.003807082510582955

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Itga2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 735.4993083333334
$ws.Range("H5").Value = 2206.497925
$ws.Range("I5").Value = 0.9467182815928301
$ws.Range("J5").Value = 0.9467182815928301
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.888921
$ws.Range("N5").Value = 5.666763
$ws.Range("O5").Value = 0.4551532417350329
$ws.Range("P5").Value = 0.4551532417350328
$ws.Range("Q5").Value = 1389.300088996308
$ws.Range("R5").Value = 12503.70080096678
$ws.Range("S5").Value = 0.4309018948767963
$ws.Range("T5").Value = 0.4309018948767963

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Itga2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 735.4993083333334
$ws.Range("H6").Value = 2206.497925
$ws.Range("I6").Value = 0.9467182815928301
$ws.Range("J6").Value = 0.9467182815928301
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.087098333333333
$ws.Range("N6").Value = 3.261295
$ws.Range("O6").Value = 0.2619465454094788
$ws.Range("P6").Value = 0.2619465454094787
$ws.Range("Q6").Value = 799.5600722569861
$ws.Range("R6").Value = 7196.040650312874
$ws.Range("S6").Value = 0.24798958333924
$ws.Range("T6").Value = 0.2479895833392399

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Itga2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 735.4993083333334
$ws.Range("H7").Value = 2206.497925
$ws.Range("I7").Value = 0.9467182815928301
$ws.Range("J7").Value = 0.9467182815928301
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.174057666666667
$ws.Range("N7").Value = 3.522173
$ws.Range("O7").Value = 0.2829002128554884
$ws.Range("P7").Value = 0.2829002128554884
$ws.Range("Q7").Value = 863.5186017767807
$ws.Range("R7").Value = 7771.667415991025
$ws.Range("S7").Value = 0.2678268033767938
$ws.Range("T7").Value = 0.2678268033767938

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col1a2"
$ws.Range("C8").Value = "Itga2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 30.939307
$ws.Range("H8").Value = 92.81792100000001
$ws.Range("I8").Value = 0.03982438491082609
$ws.Range("J8").Value = 0.03982438491082609
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.888921
$ws.Range("N8").Value = 5.666763
$ws.Range("O8").Value = 0.4551532417350329
$ws.Range("P8").Value = 0.4551532417350328
$ws.Range("Q8").Value = 58.441906717747
$ws.Range("R8").Value = 525.977160459723
$ws.Range("S8").Value = 0.01812619789226622
$ws.Range("T8").Value = 0.01812619789226622

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col1a2"
$ws.Range("C9").Value = "Itga2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 30.939307
$ws.Range("H9").Value = 92.81792100000001
$ws.Range("I9").Value = 0.03982438491082609
$ws.Range("J9").Value = 0.03982438491082609
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.087098333333333
$ws.Range("N9").Value = 3.261295
$ws.Range("O9").Value = 0.2619465454094788
$ws.Range("P9").Value = 0.2619465454094787
$ws.Range("Q9").Value = 33.63406907418833
$ws.Range("R9").Value = 302.706621667695
$ws.Range("S9").Value = 0.01043186005044827
$ws.Range("T9").Value = 0.01043186005044826

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col1a2"
$ws.Range("C10").Value = "Itga2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 30.939307
$ws.Range("H10").Value = 92.81792100000001
$ws.Range("I10").Value = 0.03982438491082609
$ws.Range("J10").Value = 0.03982438491082609
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.174057666666667
$ws.Range("N10").Value = 3.522173
$ws.Range("O10").Value = 0.2829002128554884
$ws.Range("P10").Value = 0.2829002128554884
$ws.Range("Q10").Value = 36.32453058470367
$ws.Range("R10").Value = 326.9207752623331
$ws.Range("S10").Value = 0.0112663269681116
$ws.Range("T10").Value = 0.0112663269681116
